$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'29.400.17"
$ws.Range("E2").Value = "  -0.41%  "

# Row 3
$ws.Range("D3").Value = "'1.847.56"
$ws.Range("E3").Value = "  -0.18%  "

# Row 4
$ws.Range("D4").Value = "'0.9986"

# Row 5
$ws.Range("D5").Value = "'240.33"
$ws.Range("E5").Value = "  -0.62%  "

# Row 6
$ws.Range("D6").Value = "'0.6331"
$ws.Range("E6").Value = "  +0.68%  "

# Row 7
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").Value = "'0.07563"
$ws.Range("E8").Value = "  -0.10%  "

# Row 9
$ws.Range("D9").Value = "'0.2965"
$ws.Range("E9").Value = "  -0.41%  "

# Row 10
$ws.Range("D10").Value = "'24.53"
$ws.Range("E10").Value = "  +0.80%  "

# Row 11
$ws.Range("D11").Value = "'0.07713"
$ws.Range("E11").Value = "  +0.55%  "

# Row 12
$ws.Range("D12").Value = "'1.843.43"
$ws.Range("E12").Value = "  -2.37%  "

# Row 13
$ws.Range("D13").Value = "'4.998"
$ws.Range("E13").Value = "  -0.34%  "

# Row 14
$ws.Range("D14").Value = "'0.6858"
$ws.Range("E14").Value = "  +0.13%  "

# Row 15
$ws.Range("D15").Value = "'0.00001008"
$ws.Range("E15").Value = "  +2.76%  "

# Row 16
$ws.Range("D16").Value = "'83.12"
$ws.Range("E16").Value = "  -0.84%  "

# Row 17
$ws.Range("D17").Value = "'6.167"
$ws.Range("E17").Value = "  -0.83%  "

# Row 18
$ws.Range("D18").Value = "'29.425.94"
$ws.Range("E18").Value = "  -0.49%  "

# Row 19
$ws.Range("D19").Value = "'230.20"
$ws.Range("E19").Value = "  -1.89%  "

# Row 20
$ws.Range("D20").Value = "'12.46"
$ws.Range("E20").Value = "  -0.32%  "

# Row 21
$ws.Range("D21").Value = "'1.0000"
$ws.Range("E21").Value = "  +0.06%  "

# Row 22
$ws.Range("D22").Value = "'7.577"
$ws.Range("E22").Value = "  -0.55%  "

# Row 23
$ws.Range("E23").Value = "  +0.06%  "

# Row 24
$ws.Range("D24").Value = "'156.86"
$ws.Range("E24").Value = "  +0.69%  "

# Row 25
$ws.Range("D25").Value = "'0.1400"
$ws.Range("E25").Value = "  +0.69%  "

# Row 26
$ws.Range("D26").Value = "'8.382"
$ws.Range("E26").Value = "  -0.59%  "

# Row 27
$ws.Range("D27").Value = "'17.67"
$ws.Range("E27").Value = "  -0.29%  "

# Row 28
$ws.Range("D28").Value = "'1.464"
$ws.Range("E28").Value = "  -1.20%  "

# Row 29
$ws.Range("D29").Value = "'0.05740"
$ws.Range("E29").Value = "  -1.62%  "

# Row 30
$ws.Range("D30").Value = "'1.251"
$ws.Range("E30").Value = "  -0.81%  "

# Row 31
$ws.Range("D31").Value = "'4.132"
$ws.Range("E31").Value = "  +0.54%  "

# Row 32
$ws.Range("D32").Value = "'4.027"
$ws.Range("E32").Value = "  -0.36%  "

# Row 33
$ws.Range("D33").Value = "'1.850"
$ws.Range("E33").Value = "  -2.52%  "

# Row 34
$ws.Range("D34").Value = "'1.157"
$ws.Range("E34").Value = "  -1.25%  "

# Row 35
$ws.Range("D35").Value = "'0.7172"
$ws.Range("E35").Value = "  +0.03%  "

# Row 36
$ws.Range("D36").Value = "'2.580"
$ws.Range("E36").Value = "  -0.28%  "

# Row 37
$ws.Range("D37").Value = "'1.254.40"
$ws.Range("E37").Value = "  +1.61%  "

# Row 38
$ws.Range("D38").Value = "'0.01818"
$ws.Range("E38").Value = "  +2.13%  "

# Row 39
$ws.Range("D39").Value = "'2.782"
$ws.Range("E39").Value = "  -0.80%  "

# Row 40
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'6.304"
$ws.Range("E40").Value = "  +2.78%  "

# Row 41
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.9070"
$ws.Range("E41").Value = "  -0.65%  "

# Row 42
$ws.Range("D42").Value = "'1.000"
$ws.Range("E42").Value = "  +0.09%  "

# Row 43
$ws.Range("D43").Value = "'2.004.10"
$ws.Range("E43").Value = "  -2.14%  "

# Row 44
$ws.Range("D44").Value = "'101.86"
$ws.Range("E44").Value = "  -0.01%  "

# Row 45
$ws.Range("D45").Value = "'66.25"
$ws.Range("E45").Value = "  -1.83%  "

# Row 46
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "'0.00000000118"
$ws.Range("E46").Value = "  +0.43%  "

# Row 47
$ws.Range("D47").Value = "'7.055"
$ws.Range("E47").Value = "  -3.17%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'9.154"
$ws.Range("E48").Value = "  -0.20%  "

# Row 49
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D49").Value = "'0.4032"
$ws.Range("E49").Value = "  +0.07%  "

# Row 50
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "'1.709"
$ws.Range("E50").Value = "  +1.44%  "

# Row 51
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.1132"
$ws.Range("E51").Value = "  +0.91%  "
